$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before D; this shifts existing D:K -> F:M
$ws.Columns("D:E").Insert()

$rows = @{
  7 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916)
  8 = @(6336200, 5985700, 10172600, 6508900, 12091200, 8880700, 12348700)
  9 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA")
  10 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA")
  11 = @($null, $null, $null, $null, $null, $null, $null)
  12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA")
  13 = @(0, 0, 0, 0, 0, 0, 0)
  14 = @(0, 0, 0, 0, 0, 0, 0)
  15 = @(0, 0, 0, 0, 0, 0, 0)
  16 = @($null, $null, $null, $null, $null, $null, $null)
  17 = @(5801900, 4331300, 8793600, 5020400, 13451600, 7708500, 10936900)
  18 = @(534300, 1654400, 1379000, 1488400, -1360400, 1172100, 1411800)
  19 = @($null, $null, $null, $null, $null, $null, $null)
  20 = @(0, 0, 0, 0, 0, 0, 0)
  21 = @(728600, 1786100, 1505500, 1591900, -1245100, 1276300, 1511500)
  22 = @(267200, 232200, 236700, 212800, 219500, 227700, 207600)
  23 = @(267200, 1422200, 1142400, 1275600, -1580000, 944400, 1204100)
  24 = @(177900, 96700, 183100, 250800, -666800, 9700, 226200)
  25 = @(0, 0, 0, 0, 0, 0, 0)
  26 = @(89300, 1325400, 959300, 1024800, -913100, 934700, 977900)
  27 = @(263500, 1047100, 906400, 992000, -873700, 793300, 905000)
  28 = @(0, 0, 0, 0, 0, 0, 0)
  29 = @(145900, 92300, "NA", "NA", -351300, "NA", "NA")
  30 = @(0, 0, 0, 0, 0, 0, 0)
  31 = @(0, 0, 0, 0, 0, 0, 0)
  32 = @(0, 0, 0, 0, 0, 0, 0)
  33 = @(409300, 1139400, 906400, 992000, -1225000, 793300, 905000)
  34 = @(0, 0, 0, 0, 0, 0, 0)
  35 = @(409300, 1139400, 906400, 992000, -1225000, 793300, 905000)
  38 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916)
  39 = @($null, $null, $null, $null, $null, $null, $null)
  40 = @($null, $null, $null, $null, $null, $null, $null)
  41 = @(12067400, 11640900, 12822700, 12666500, 11881300, 11349900, 11807600)
  42 = @(0, 0, 0, 0, 0, 0, 0)
  43 = @(0, 0, 0, 0, 0, 0, 0)
  44 = @(0, 0, 0, 0, 0, 0, 0)
  45 = @(0, 0, 0, 0, 0, 0, 0)
  46 = @(0, 0, 0, 0, 0, 0, 0)
  47 = @(224524400, 218567000, 221062400, 217306300, 212937100, 207502100, 210094200)
  48 = @(0, 0, 0, 0, 0, 0, 0)
  49 = @(7514300, 7421300, 7490500, 7427200, 7323000, 7284300, 7440600)
  50 = @(0, 0, 0, 0, 0, 0, 0)
  51 = @(0, 0, 0, 0, 0, 0, 0)
  52 = @(3213500, 3193400, 3553600, 3570000, 3400300, 3464300, 3389900)
  53 = @(0, 0, 0, 0, 0, 0, 0)
  54 = @(558359200, 556548500, 559813400, 550583700, 542925800, 530901600, 540107400)
  55 = @($null, $null, $null, $null, $null, $null, $null)
  56 = @($null, $null, $null, $null, $null, $null, $null)
  57 = @(0, 0, 0, 0, 0, 0, 0)
  58 = @(0, 0, 0, 0, 0, 0, 0)
  59 = @(480110700, 480207500, 484234400, 476036200, 470369000, 456448600, 463866800)
  60 = @(0, 0, 0, 0, 0, 0, 0)
  61 = @(10047600, 9948600, 10040100, 9809400, 9802700, 9871900, 9802000)
  62 = @(1350000, 1120800, 1084300, 1077600, 953300, 1400600, 1481700)
  63 = @(0, 0, 0, 0, 0, 0, 0)
  64 = @(0, 0, 0, 0, 0, 0, 0)
  65 = @(0, 0, 0, 0, 0, 0, 0)
  66 = @(524082400, 523529400, 526862700, 518503000, 512239000, 498995800, 507829600)
  67 = @($null, $null, $null, $null, $null, $null, $null)
  68 = @(0, 0, 0, 0, 0, 0, 0)
  69 = @(0, 0, 0, 0, 0, 0, 0)
  70 = @(2844400, 2844400, 2844400, 2844400, 2662000, 2662000, 2662000)
  71 = @(0, 0, 0, 0, 0, 0, 0)
  72 = @(9454400, 9572000, 8757900, 8176600, 7503900, 9031000, 8539800)
  73 = @(0, 0, 0, 0, 0, 0, 0)
  74 = @(0, 0, 0, 0, 0, 0, 0)
  75 = @(0, 0, 0, 0, 0, 0, 0)
  76 = @(31432500, 30174700, 30106300, 29236300, 28024700, 29243700, 29615800)
  77 = @(0, 0, 0, 0, 0, 0, 0)
  80 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916)
  81 = @(409300, 1139400, 906400, 992000, -1225000, 793300, 905000)
  82 = @($null, $null, $null, $null, $null, $null, $null)
  83 = @(194200, 131700, 126500, 103400, 115400, 104200, 99700)
  84 = @(0, 0, 0, 0, 0, 0, 0)
  85 = @(0, 0, 0, 0, 0, 0, 0)
  86 = @(0, 0, 0, 0, 0, 0, 0)
  87 = @(0, 0, 0, 0, 0, 0, 0)
  88 = @(0, 0, 0, 0, 0, 0, 0)
  89 = @(3524600, 4199600, 3692000, 2863700, 3822300, 3422600, 3625000)
  90 = @($null, $null, $null, $null, $null, $null, $null)
  91 = @(0, 0, 0, 0, 0, 0, 0)
  92 = @(0, 0, 0, 0, 0, 0, 0)
  93 = @(0, 0, 0, 0, 0, 0, 0)
  94 = @(-2662000, -5257800, -3439000, -2782600, -2972400, -3754500, -2799700)
  95 = @($null, $null, $null, $null, $null, $null, $null)
  96 = @(-267900, -356500, -358000, -348300, -331900, -331200, -331200)
  97 = @(0, 0, 0, 0, 0, 0, 0)
  98 = @(0, 0, 0, 0, 0, 0, 0)
  99 = @(0, 0, 0, 0, 0, 0, 0)
  100 = @(-772500, -21600, -173400, 428700, -456900, 151800, -649000)
  101 = @(425700, -175600, 104900, 256800, 44700, -303600, -212100)
  102 = @(515700, -1255500, 184600, 766500, 437600, -483700, -35700)
}

$dateRows = @(7, 38, 80)

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $rng = $ws.Range("D" + $r + ":J" + $r)
  if ($dateRows -contains $r) {
    $rng.NumberFormat = "[$-409]d\-mmm\-yy;@"
  } else {
    $rng.NumberFormat = "#,##0"
  }
  for ($i = 0; $i -lt 7; $i++) {
    $cell = $ws.Cells.Item($r, 4 + $i)
    $v = $vals[$i]
    if ($null -eq $v) {
      $cell.Value = $null
    } else {
      $cell.Value = $v
    }
  }
}

$ws.Columns("D:E").ColumnWidth = 15.083333333333334